$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": append a new day row (row 51) below the existing data (row 50)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(51).Insert()
$wsAll.Range("A51").Value = 43979
$wsAll.Range("B51").Value = 285
$wsAll.Range("C51").Value = 282
$wsAll.Range("D51").Value = 15
$wsAll.Range("E51").Value = 12
$wsAll.Range("F51").Value = 3
$wsAll.Range("G51").Value = 12
$wsAll.Range("H51").Value = 255

$wsAll.Activate()
$wsAll.Range("B51:H51").Select()

# ---------------------------------------------------------------------------
# Sheet "kobe": insert the new day row (106) just above the footnote row,
# which shifts the footnote row from 106 down to 107
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(106).Insert()
$wsKobe.Range("A106").Value = 43979
$wsKobe.Range("B106").Value = 15
$wsKobe.Range("C106").Value = 3121
$wsKobe.Range("E106").Value = 285
$wsKobe.Range("F106").Value = 12
$wsKobe.Range("G106").Value = 10
$wsKobe.Range("H106").Value = 2
$wsKobe.Range("I106").Value = 12
$wsKobe.Range("J106").Value = 244

$wsKobe.Activate()
$wsKobe.Range("B106:J106").Select()

# ---------------------------------------------------------------------------
# Sheet "other": insert the new day row (81) just above the footnote row,
# which shifts the footnote row from 81 down to 82
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(81).Insert()
$wsOther.Range("A81").Value = 43979
$wsOther.Range("B81").Value = 0
$wsOther.Range("C81").Value = 14
$wsOther.Range("D81").Value = 3
$wsOther.Range("E81").Value = 2
$wsOther.Range("F81").Value = 1
$wsOther.Range("G81").Value = 0
$wsOther.Range("H81").Value = 11

# "other" becomes the active / selected sheet after the update
$wsOther.Activate()
$wsOther.Range("B85").Select()

Write-Host "edit complete"
